$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.924.05'
$ws.Range('E2').Value = '  -0.53%  '
$ws.Range('D3').Value = '3.033.19'
$ws.Range('E3').Value = '  -0.67%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '586.97'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.05%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '149.21'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.67%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.527'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.15%  '
$ws.Range('D9').Value = '3.031.28'
$ws.Range('E9').Value = '  -0.75%  '
$ws.Range('E10').Value = '  -2.90%  '
$ws.Range('E11').Value = '  +0.02%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.445'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.22%  '
$ws.Range('E13').Value = '  -1.97%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.35'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -4.35%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.121'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.40%  '
$ws.Range('D16').Value = '3.537.63'
$ws.Range('E16').Value = '  -0.60%  '
$ws.Range('E17').Value = '  -0.39%  '
$ws.Range('D18').Value = '62.844.28'
$ws.Range('E18').Value = '  -0.73%  '
$ws.Range('D19').Value = '3.032.76'
$ws.Range('E19').Value = '  -0.65%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '469.94'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.38%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.04'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.693'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.60%  '
$ws.Range('E23').Value = '  -0.96%  '
$ws.Range('E24').Value = '  -1.81%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '80.93'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.51%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.45'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.51%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.46'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.34%  '
$ws.Range('E28').Value = '  +0.26%  '
$ws.Range('E29').Value = '  +0.05%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.25'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.68%  '
$ws.Range('E31').Value = '  -0.57%  '
$ws.Range('E32').Value = '  +0.45%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '27.78'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.76%  '
$ws.Range('E34').Value = '  -4.13%  '
$ws.Range('E35').Value = '  +0.11%  '
$ws.Range('D36').Value = '0.0₃0807'
$ws.Range('E36').Value = '  -1.52%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.80'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.75%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '50.39'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.54%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '9.03'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.31%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.96'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -10.33%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '427.11'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.11%  '
$ws.Range('E43').Value = '  +2.62%  '
$ws.Range('E44').Value = '  -2.40%  '
$ws.Range('D45').Value = '2.807.20'
$ws.Range('E45').Value = '  +0.74%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0357'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.21%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '37.81'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -8.98%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '129.12'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.75%  '
$ws.Range('E49').Value = '  -0.02%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '24.48'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.19%  '
$ws.Range('E51').Value = '  -0.09%  '
